$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText, $matchWholeWord) {
    $d.Content.Find.Execute($findText, $true, $matchWholeWord, $false, $false, $false, `
                             $true, 1, $false, $replaceText, 2)
}

# 1. Фамилия: Семенов -> dfghsdf
Replace-Text "Семенов" "dfghsdf" $true

# 2. Имя: Андрей -> (empty)
Replace-Text "Андрей" "" $true

# 3. Отчество: Игоревич -> (empty)  (only the standalone occurrence)
Replace-Text "Игоревич" "" $true

# 4. Дата рождения: 25 мая 1979 -> 01 января 1970
Replace-Text "25 мая 1979" "01 января 1970" $false

# 5. Место жительства: МО, Красногорск, Светлая, 11,  -> (empty)
Replace-Text "МО, Красногорск, Светлая, 11, " "" $false

# 6. Место работы -> (empty)
Replace-Text 'Филиал №1 ФГБУ "НМИЦ ВМТ им.А.А.Вишневского" МО РФ' "" $false

# 7. Удостоверение личности -> (empty)
Replace-Text "Паспорт РФ: 1817 314592, выдан: отделом по вопросам миграции отдела полиции №1 УМВД России по г.Волгограду 2017-05-10" "" $false

# 8. Контактный телефон -> (empty)
Replace-Text "8-905-396-65-40" "" $false

# 9. Full name occurrences ("Я Семенов Андрей Игоревич") -> "dfghsdf  " (x2)
Replace-Text "Семенов Андрей Игоревич" "dfghsdf  " $false

# 10. Blood group: B(III) -> 0(I)
Replace-Text "B(III)" "0(I)" $false

# 11. Rh factor: Rh+ -> Rh-
Replace-Text "Rh+" "Rh-" $false

# 12. Phenotype: D+C-E-c-e-K- -> D-C-E-c-e-K-
Replace-Text "D+C-E-c-e-K-" "D-C-E-c-e-K-" $false
